# Build site at 2023-04-12 14:53:07 UTC
# Update LOM3258.xlsx course-plan sheet: fix several mis-shifted cell values,
# insert the "teacher name" row under "Docentes responsaveis:", and fill in
# the previously-empty Objetivos / Programa resumido / Programa / Bibliografia
# (and the off-by-one Metodo/Criterio/Norma de recuperacao) texts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Fix "Objetivos:" row (10) - was showing the teacher's name by mistake
# ------------------------------------------------------------------
$ws.Range("B10").Value = "Proporcionar ao aluno ingressante de Engenharia Física os conhecimentos práticos de eletrônica e computação física com microcontrolador Arduino visando sua aplicação em projetos científicos e tecnológicos."
$ws.Range("C10").Value = "Proporcionar ao aluno ingressante de Engenharia Física os conhecimentos práticos de eletrônica e computação física com microcontrolador Arduino visando sua aplicação em projetos científicos e tecnológicos."

# ------------------------------------------------------------------
# 2) Insert a new row 13 under "Docentes responsaveis:" (row 12) to hold
#    the teacher's name, pushing everything from the old row 13 down by one.
# ------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# Copy B11:C11 formatting (style 2 / style 3, no row height) onto the new
# row 13 cells, then set their values and clear the stray A13 formatting
# that Insert() carried down from the row above.
$ws.Range("B11:C11").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("A13").Clear()

# ------------------------------------------------------------------
# 3) Fix "Programa resumido:" row (now row 14) - was "Semestral"
# ------------------------------------------------------------------
$ws.Range("B14").Value = "Introdução ao Arduino. Conceitos de eletrônica analógica e digital. Montagem de circuitos eletrônicos básicos. Programação e controle de circuitos eletrônicos em linguagem C. Aplicação e desenvolvimento de projetos baseados em Arduino."
$ws.Range("C14").Value = "Introdução ao Arduino. Conceitos de eletrônica analógica e digital. Montagem de circuitos eletrônicos básicos. Programação e controle de circuitos eletrônicos em linguagem C. Aplicação e desenvolvimento de projetos baseados em Arduino."

# ------------------------------------------------------------------
# 4) Fix "Programa:" row (now row 16) - was "01/01/2023"
# ------------------------------------------------------------------
$ws.Range("B16").Value = "Introdução ao microcontrolador Arduino: histórico, tipos e recursos. Oficina prática: instalação e configuração do IDE Arduino.Conceitos básicos de eletrônica: funcionamento da protoboard, componentes e instrumentos eletrônicos, medições com multímetro e osciloscópio. Grandezas elétricas: resistência, tensão e corrente. Oficina: montagem de circuitos eletrônicos.Introdução à linguagem de programação Wiring baseada em C/C++. Tipos de dados, sintaxe básica, controle de fluxo, funções da biblioteca padrão. Principais bibliotecasEntradas e saídas do Arduino. Sinais analógicos e digitais.Controle de dispositivos utilizando PWM.Eletrônica analógica. Conversores analógico-digitais do Arduino. Oficina: leitura de dados de sensores. Comunicação serial/USB com o PC. Utilização do Monitor Serial da IDE.Controle de motor cc e servomotor com PWM. Controle de potência com relé e SSR.Tópicos avançados: comunicação Ethernet com Arduino. Comunicação sem fio via Bluetooth.Armazenamento de dados utilizando a EEPROM do ATMega328 e cartão de memória SD.Desenvolvimento de software de qualidade.Desenvolvimento de projetos utilizando microcontrolador Arduino."
$ws.Range("C16").Value = "Introdução ao microcontrolador Arduino: histórico, tipos e recursos. Oficina prática: instalação e configuração do IDE Arduino.Conceitos básicos de eletrônica: funcionamento da protoboard, componentes e instrumentos eletrônicos, medições com multímetro e osciloscópio. Grandezas elétricas: resistência, tensão e corrente. Oficina: montagem de circuitos eletrônicos.Introdução à linguagem de programação Wiring baseada em C/C++. Tipos de dados, sintaxe básica, controle de fluxo, funções da biblioteca padrão. Principais bibliotecasEntradas e saídas do Arduino. Sinais analógicos e digitais.Controle de dispositivos utilizando PWM.Eletrônica analógica. Conversores analógico-digitais do Arduino. Oficina: leitura de dados de sensores. Comunicação serial/USB com o PC. Utilização do Monitor Serial da IDE.Controle de motor cc e servomotor com PWM. Controle de potência com relé e SSR.Tópicos avançados: comunicação Ethernet com Arduino. Comunicação sem fio via Bluetooth.Armazenamento de dados utilizando a EEPROM do ATMega328 e cartão de memória SD.Desenvolvimento de software de qualidade.Desenvolvimento de projetos utilizando microcontrolador Arduino."

# ------------------------------------------------------------------
# 5) Fix the off-by-one Avaliacao block (rows 19-22 after the insert):
#    Metodo had the teacher's name, Criterio had Metodo's text, Norma de
#    recuperacao had Criterio's text, and Bibliografia had Norma's text.
# ------------------------------------------------------------------
$ws.Range("B19").Value = "Aulas expositivas, práticas e de realização de projetos."
$ws.Range("C19").Value = "Aulas expositivas, práticas e de realização de projetos."

$ws.Range("B20").Value = "Média das notas de trabalhos, atividades e relatório de projeto."
$ws.Range("C20").Value = "Média das notas de trabalhos, atividades e relatório de projeto."

$ws.Range("B21").Value = "Devido às características da disciplina não será oferecida recuperação."
$ws.Range("C21").Value = "Devido às características da disciplina não será oferecida recuperação."

$ws.Range("B22").Value = "BANZI, M. Primeiros passos com o Arduino, São Paulo: O´Reilly Novatec, 2010.`nMcROBERTS, M. Arduino Básico, São Paulo: Novatec, 2011.`nMONK, S. Programação com Arduino, Porto Alegre: Bookman Editora, 2013.`nMONK, S. Programação com Arduino II, Porto Alegre: Bookman Editora, 2015.`nBLUM, J. Exploring Arduino, New York: John Wiley, 2013."
$ws.Range("C22").Value = "BANZI, M. Primeiros passos com o Arduino, São Paulo: O´Reilly Novatec, 2010.`nMcROBERTS, M. Arduino Básico, São Paulo: Novatec, 2011.`nMONK, S. Programação com Arduino, Porto Alegre: Bookman Editora, 2013.`nMONK, S. Programação com Arduino II, Porto Alegre: Bookman Editora, 2015.`nBLUM, J. Exploring Arduino, New York: John Wiley, 2013."

# ------------------------------------------------------------------
# 6) Column A was sized together with column B (min=1,max=2); column B
#    already carries its own 60.7109375 override, so narrow column A's
#    band down to itself so it no longer spans column B.
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30.7109375

Write-Host "edit complete"
